$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "28.465.99"
$dCell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +2.15%  "
$dCell = $ws.Cells.Item(3, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.827.36"
$dCell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +1.89%  "
$ws.Cells.Item(4, 5).Value = "  +0.18%  "
$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "315.52"
$dCell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.34%  "
$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.002"
$dCell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.16%  "
$dCell = $ws.Cells.Item(7, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.5050"
$dCell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -5.04%  "
$dCell = $ws.Cells.Item(8, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.3912"
$dCell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +1.20%  "
$dCell = $ws.Cells.Item(9, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.07646"
$dCell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +2.66%  "
$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "41.90"
$dCell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +1.17%  "
$dCell = $ws.Cells.Item(11, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.104"
$dCell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.73%  "
$dCell = $ws.Cells.Item(12, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "21.03"
$dCell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +3.46%  "
$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.284"
$dCell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.63%  "
$ws.Cells.Item(14, 5).Value = "  +0.22%  "
$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "7.562"
$dCell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.60%  "
$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.828.05"
$dCell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.86%  "
$dCell = $ws.Cells.Item(17, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "93.22"
$dCell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +5.57%  "
$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.00001084"
$dCell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +2.37%  "
$dCell = $ws.Cells.Item(19, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.06651"
$dCell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.79%  "
$dCell = $ws.Cells.Item(20, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "17.70"
$dCell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.70%  "
$ws.Cells.Item(21, 5).Value = "  +0.08%  "
$dCell = $ws.Cells.Item(22, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.162"
$dCell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +3.24%  "
$dCell = $ws.Cells.Item(23, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "28.495.40"
$dCell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +2.12%  "
$dCell = $ws.Cells.Item(24, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "11.12"
$dCell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.05%  "
$dCell = $ws.Cells.Item(25, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.258"
$dCell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +7.95%  "
$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "156.70"
$dCell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.27%  "
$dCell = $ws.Cells.Item(27, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "20.59"
$dCell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +2.24%  "
$dCell = $ws.Cells.Item(28, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.035.88"
$dCell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.72%  "
$dCell = $ws.Cells.Item(29, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.398"
$dCell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +4.07%  "
$dCell = $ws.Cells.Item(30, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "125.03"
$dCell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +2.81%  "
$dCell = $ws.Cells.Item(31, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.127"
$dCell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +2.67%  "
$dCell = $ws.Cells.Item(32, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.1082"
$dCell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.87%  "
$dCell = $ws.Cells.Item(33, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "5.673"
$dCell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +3.08%  "
$dCell = $ws.Cells.Item(34, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.662"
$dCell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.14%  "
$dCell = $ws.Cells.Item(35, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.07061"
$dCell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +1.40%  "
$dCell = $ws.Cells.Item(36, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.2221"
$dCell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.90%  "
$dCell = $ws.Cells.Item(37, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "8.930"
$dCell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +6.26%  "
$dCell = $ws.Cells.Item(38, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.02320"
$dCell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +2.18%  "
$dCell = $ws.Cells.Item(39, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "5.136"
$dCell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.54%  "
$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.6240"
$dCell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +2.01%  "
$dCell = $ws.Cells.Item(41, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "11.21"
$dCell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.41%  "
$dCell = $ws.Cells.Item(42, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.182"
$dCell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.53%  "
$ws.Cells.Item(43, 5).Value = "  +0.11%  "
$dCell = $ws.Cells.Item(44, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.399"
$dCell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -1.08%  "
$ws.Cells.Item(45, 5).Value = "  +1.22%  "
$dCell = $ws.Cells.Item(46, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.5897"
$dCell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +3.29%  "
$dCell = $ws.Cells.Item(47, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.720"
$dCell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.09%  "
$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "124.34"
$dCell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.48%  "
$dCell = $ws.Cells.Item(49, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.979"
$dCell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.45%  "
$dCell = $ws.Cells.Item(50, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.193"
$dCell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.45%  "
$dCell = $ws.Cells.Item(51, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.06912"
$dCell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.70%  "
